# Applies the "agregando fecha de descripción" change:
# Adds a note "dia/mes/año Ejemplo: 11/09/2001" to the "Observaciones" column
# for the "Fecha de nacimiento" (row 14) and "Fecha de ingreso" (row 19) rows
# on the "Descripción por campo" sheet, sets C14 to 0, and makes that sheet
# the active tab with the selection left on D19 (topLeftCell A13).

$wb = $excel.ActiveWorkbook

$wsProp = $wb.Worksheets.Item("Formato propuesto")
$wsDesc = $wb.Worksheets.Item("Descripción por campo")

# Fill in the new data on the "Descripción por campo" sheet.
$wsDesc.Range("C14").Value = 0
$wsDesc.Range("D14").Value = "dia/mes/año Ejemplo: 11/09/2001"
$wsDesc.Range("D19").Value = "dia/mes/año Ejemplo: 11/09/2001"

# Make "Descripción por campo" the active sheet/tab, with the view scrolled
# so row 13 is at the top and D19 selected.
$wsDesc.Activate()
$wsDesc.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 13

$wb.Save()
